# "retrait fichiers inutiles + finalisation mail + logo + modifs test excel"
#
# This script:
#   - adds the client e-mail address (with a mailto hyperlink) for rows 17
#     and 18 in column C
#   - duplicates row 17 into a brand-new row 19 (same data, new GlobalID)
#     which also gets its e-mail hyperlink in column C
#   - updates the saved window selection (D23, no frozen/scrolled topLeftCell)
#   - leaves every other existing cell untouched

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-MailHyperlink {
    param($cellRef, $email)

    $ws.Range($cellRef).Value = $email
    $ws.Hyperlinks.Add($ws.Range($cellRef), "mailto:" + $email)

    # `Hyperlinks.Add` auto-applies Excel's built-in blue/underlined
    # "Hyperlink" style - put the cell back on the sheet's normal
    # (green-filled) look used by the rest of the row.
    $ws.Range($cellRef).Font.Underline = $false
    $ws.Range($cellRef).Font.Name = "Calibri"
    $ws.Range($cellRef).Font.Size = 11
    $ws.Range($cellRef).Font.Color = 0
    $ws.Range($cellRef).Interior.Color = $ws.Range("B17").Interior.Color()
}

# --- row 17 : quochiepdao92@gmail.com -------------------------------------
Set-MailHyperlink "C17" "quochiepdao92@gmail.com"

# --- row 18 : lothesven@yahoo.fr -------------------------------------------
Set-MailHyperlink "C18" "lothesven@yahoo.fr"

# --- new row 19 : copy of row 17, new GlobalID ------------------------------
$ws.Range("A19").Value = "60222876-f542-4b59-a34d-e2ecc13c45d9"
$ws.Range("B19").Value = $ws.Range("B17").Value()
$ws.Range("D19").Value = $ws.Range("D17").Value()
$ws.Range("E19").Value = $ws.Range("E17").Value()
$ws.Range("F19").Value = $ws.Range("F17").Value()
$ws.Range("F19").NumberFormat = $ws.Range("F17").NumberFormat()
$ws.Range("G19").Value = $ws.Range("G17").Value()
$ws.Range("H19").Value = $ws.Range("H17").Value()
$ws.Range("I19").Value = $ws.Range("I17").Value()
$ws.Range("I19").NumberFormat = $ws.Range("I17").NumberFormat()
$ws.Range("J19").Value = $ws.Range("J17").Value()
$ws.Range("L19").Value = $ws.Range("L17").Value()
$ws.Range("N19").Value = $ws.Range("N17").Value()
$ws.Range("O19").Value = $ws.Range("O17").Value()

Set-MailHyperlink "C19" "quochiepdao92@gmail.com"

# --- window selection : D23, scrolled back to column A ---------------------
$ws.Range("D23").Select()

Write-Output "done"
